# Reorder sheets: move "Observables" and "Functions" to just after "Concentrations",
# and move "Stop conditions" to just after "Parameters".
$wb = $excel.ActiveWorkbook

$obs = $wb.Worksheets.Item("Observables")
$concentrations = $wb.Worksheets.Item("Concentrations")
$obs.Move($null, $concentrations)

$func = $wb.Worksheets.Item("Functions")
$obsAfterMove = $wb.Worksheets.Item("Observables")
$func.Move($null, $obsAfterMove)

$stop = $wb.Worksheets.Item("Stop conditions")
$params = $wb.Worksheets.Item("Parameters")
$stop.Move($null, $params)

# Add header rows for the new model sheets.
$obsSheet = $wb.Worksheets.Item("Observables")
$obsSheet.Range("A1").Value = "Id"
$obsSheet.Range("B1").Value = "Name"
$obsSheet.Range("C1").Value = "Model"
$obsSheet.Range("D1").Value = "Species"
$obsSheet.Range("E1").Value = "Observables"
$obsSheet.Range("F1").Value = "Comments"
$obsSheet.Range("A1:F1").Select() | Out-Null

$funcSheet = $wb.Worksheets.Item("Functions")
$funcSheet.Range("A1").Value = "Id"
$funcSheet.Range("B1").Value = "Name"
$funcSheet.Range("C1").Value = "Model"
$funcSheet.Range("D1").Value = "Expression"
$funcSheet.Range("E1").Value = "Comments"
$funcSheet.Range("A1:E1").Select() | Out-Null

$stopSheet = $wb.Worksheets.Item("Stop conditions")
$stopSheet.Range("A1").Value = "Id"
$stopSheet.Range("B1").Value = "Name"
$stopSheet.Range("C1").Value = "Model"
$stopSheet.Range("D1").Value = "Expression"
$stopSheet.Range("E1").Value = "Comments"
$stopSheet.Range("A2").Select() | Out-Null

# Make "Stop conditions" the active/selected sheet (matches new activeTab index).
$stopSheet.Activate()
